# This workbook's data rows (2-21) were reordered (rows rearranged into a
# new order); no cell content was otherwise changed. Column headers stay in
# row 1. We reproduce the new row order by reading each source row's full
# values into memory, then writing them back in the new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 21
$firstCol = 1   # A
$lastCol = 51   # AY

# Mapping: new row number (2..21) -> original row number it should contain.
$rowMap = @{
    2  = 4
    3  = 5
    4  = 9
    5  = 15
    6  = 16
    7  = 17
    8  = 18
    9  = 19
    10 = 21
    11 = 2
    12 = 3
    13 = 6
    14 = 7
    15 = 8
    16 = 10
    17 = 11
    18 = 12
    19 = 13
    20 = 14
    21 = 20
}

# Columns whose values must stay plain text and not be silently
# reinterpreted (as a date serial number or as a number) when written
# back. Y/Z/AA/AB hold date/time text like "2018-05-14" or "00:00", and
# I ("Antal") holds numeric-looking counts that are nonetheless stored
# as text in this workbook.
$textColumns = @(9, 25, 26, 27, 28)   # I, Y, Z, AA, AB

# Snapshot the original values of every data row before we start
# overwriting anything.
$originalRows = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $rng = $ws.Range($ws.Cells.Item($r, $firstCol), $ws.Cells.Item($r, $lastCol))
    $originalRows[$r] = $rng.Value2
}

# Make sure the date/time columns are formatted as Text so that assigning
# strings such as "2018-05-14" does not get silently converted into a
# date serial number.
foreach ($col in $textColumns) {
    $colRng = $ws.Range($ws.Cells.Item($firstDataRow, $col), $ws.Cells.Item($lastDataRow, $col))
    $colRng.NumberFormat = "@"
}

# Write each new row using the values captured from its source row.
for ($newRow = $firstDataRow; $newRow -le $lastDataRow; $newRow++) {
    $srcRow = $rowMap[$newRow]
    $srcVals = $originalRows[$srcRow]

    $destRng = $ws.Range($ws.Cells.Item($newRow, $firstCol), $ws.Cells.Item($newRow, $lastCol))
    $destRng.Value2 = $srcVals
}

Write-Output "Reordered rows $firstDataRow-$lastDataRow"
